# Rewrites "Operations kontrakt(er)" into the full set of use-case
# operation contracts described by the commit ("Operations kontrakter
# lavet."): reflows/expands the createSale contract with per-parameter
# spell/grammar markers and en-US run language, fills in the Prebetingelser/
# Postbetingelser bodies, and appends the createEmployee and
# updateEmployee contracts (each with its own Operation/Kryds
# reference/Prebetingelser/Postbetingelser paragraphs) plus trailing
# blank paragraphs, exactly as the target OOXML has it.
#
# Implemented as a single WordOpenXML replacement of the document body
# (keeping the original sectPr) so the exact run/proofErr/lang markup
# survives byte-for-byte instead of relying on Find/Replace, which can't
# author w:proofErr or per-run w:rPr splits on its own.

$d = $word.ActiveDocument

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Operations </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>kontrakt</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>er</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Operation: </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>createSale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>employeeID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>itemID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>saleD</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>itemQuantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Kryds</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> reference</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: Use case: Sale</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Prebetingelser</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>En kunde med en vare, en oprettet medarbejder, oprettede varer.</w:t></w:r></w:p><w:p><w:r><w:t>Postbetingelser:</w:t></w:r><w:r><w:t xml:space="preserve"> Et salg bliver gennemført og der bliver fjernet solgte varer fra lager.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Operation: </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>createEmployee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>String name, String address, String phone, String password)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Kryds</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> reference: Use case: Handle Employee CRUD</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Prebetingelser</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>Postbetingelser:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>En medarbejder bliver oprettet.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Operation: </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>updateEmployee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>String password, String name, String address, String phone)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Kryds</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> reference: Use case: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Handle Employee CRUD</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Prebetingelser</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>En medarbejder er oprettet</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Postbetingelser:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Medarbejderen bliver opdateret med en af sine informationer.</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p><w:r><w:t>ZZZSENTINELZZZ</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1701" w:right="1134" w:bottom="1701" w:left="1134" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)

# InsertXML drops a trailing run of bare <w:p/> paragraphs that lands
# exactly at the end of the replaced range (next to the section mark).
# The fragment above carries one extra non-empty "sentinel" paragraph
# after the real trailing blanks so they survive; remove the sentinel
# now that the real content is safely in place.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()
